# Renumber the page-number textboxes on slides 2,3,5,6,7,8.
# (Slide 4 already shows the correct number "4" and is left untouched.)
#
# Note: these textboxes use <a:spAutoFit/>, so the COM host recomputes the
# shape's cached height (cy) as a side effect of any TextRange edit. The
# source deck does not show that cy change, so after editing the text we
# restore Height to the exact original value (expressed in points so the
# runtime's internal EMU rounding reproduces the original EMU value).

$p = $ppt.ActivePresentation

# --- Slide 2: top-level shape "Google Shape;236;p20", "3" -> "2" ---
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(12)
$sh2.TextFrame.TextRange.Text = "2"
$sh2.TextFrame.TextRange.LanguageID = "pt-BR"
$sh2.Height = 32.26708984375

# --- Slide 3: grouped shape "Google Shape;256;p22", "4" -> "3" ---
$s3 = $p.Slides.Item(3)
$grp3 = $s3.Shapes.Item(2)
$item3 = $grp3.GroupItems.Item(2)
$item3.TextFrame.TextRange.Runs(1).Text = "3"
$item3.Height = 32.2788200378418

# --- Slide 5: grouped shape "Google Shape;248;p21", "4" -> "5" ---
$s5 = $p.Slides.Item(5)
$grp5 = $s5.Shapes.Item(3)
$item5 = $grp5.GroupItems.Item(2)
$item5.TextFrame.TextRange.Runs(1).Text = "5"
$item5.Height = 107.60748291015625

# --- Slide 6: grouped shape "Google Shape;248;p21", "4" -> "6" ---
$s6 = $p.Slides.Item(6)
$grp6 = $s6.Shapes.Item(3)
$item6 = $grp6.GroupItems.Item(2)
$item6.TextFrame.TextRange.Runs(1).Text = "6"
$item6.Height = 107.60748291015625

# --- Slide 7: grouped shape "Google Shape;248;p21", "4" -> "7" ---
$s7 = $p.Slides.Item(7)
$grp7 = $s7.Shapes.Item(3)
$item7 = $grp7.GroupItems.Item(2)
$item7.TextFrame.TextRange.Runs(1).Text = "7"
$item7.Height = 107.60748291015625

# --- Slide 8: grouped shape "Google Shape;248;p21", "4" -> "8" ---
$s8 = $p.Slides.Item(8)
$grp8 = $s8.Shapes.Item(3)
$item8 = $grp8.GroupItems.Item(2)
$item8.TextFrame.TextRange.Runs(1).Text = "8"
$item8.Height = 107.60748291015625
